$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header/value labels (order matters for shared-string table ordering)
$ws.Range("A1").Value = "TestCase"
$ws.Range("A2").Value = "TempVariance"
$ws.Range("B1").Value = "Temperature(in Celsius)"
$ws.Range("B2").Value = 29.26

# Remove the now-unused columns C and D
$ws.Range("C1:D2").ClearContents()

# Set explicit column widths (A, B, C) matching the authored layout
$ws.Columns.Item(1).ColumnWidth = 10.666666666666666
$ws.Columns.Item(2).ColumnWidth = 18.833333333333336
$ws.Columns.Item(3).ColumnWidth = 16.833333333333336

# Select B2 to match the saved view state
$ws.Range("B2").Select()
